$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.389.29"
$ws.Range("E2").Value = "  -4.13%  "
$ws.Range("D3").Value = "1.564.64"
$ws.Range("E3").Value = "  -4.06%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "290.54"
$ws.Range("E6").Value = "  -2.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3687"
$ws.Range("E7").Value = "  -2.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.28"
$ws.Range("E8").Value = "  -1.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3398"
$ws.Range("E9").Value = "  -2.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.166"
$ws.Range("E10").Value = "  -3.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07649"
$ws.Range("E11").Value = "  -4.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.36"
$ws.Range("E13").Value = "  -2.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.052"
$ws.Range("E14").Value = "  -3.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.906"
$ws.Range("E15").Value = "  -4.38%  "
$ws.Range("D16").Value = "1.566.35"
$ws.Range("E16").Value = "  -3.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001125"
$ws.Range("E17").Value = "  -5.78%  "
$ws.Range("E18").Value = "  -4.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06734"
$ws.Range("E19").Value = "  -2.99%  "
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.249"
$ws.Range("E21").Value = "  -5.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.53"
$ws.Range("E22").Value = "  -4.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5302"
$ws.Range("E23").Value = "  -6.85%  "
$ws.Range("E24").Value = "  -2.91%  "
$ws.Range("D25").Value = "22.406.13"
$ws.Range("E25").Value = "  -4.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.354"
$ws.Range("E26").Value = "  -2.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.838"
$ws.Range("E27").Value = "  -3.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.05"
$ws.Range("E28").Value = "  -4.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "146.64"
$ws.Range("E29").Value = "  -1.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.978"
$ws.Range("E30").Value = "  -3.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "125.39"
$ws.Range("E31").Value = "  -4.27%  "
$ws.Range("D32").Value = "1.737.01"
$ws.Range("E32").Value = "  -4.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.018"
$ws.Range("E33").Value = "  +3.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.186"
$ws.Range("E34").Value = "  -8.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.014"
$ws.Range("E35").Value = "  -4.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.09"
$ws.Range("E36").Value = "  -9.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08527"
$ws.Range("E37").Value = "  -2.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02540"
$ws.Range("E38").Value = "  -4.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2319"
$ws.Range("E39").Value = "  -3.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.511"
$ws.Range("E40").Value = "  -5.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06451"
$ws.Range("E41").Value = "  -5.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.274"
$ws.Range("E42").Value = "  -1.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.67"
$ws.Range("E43").Value = "  -8.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6341"
$ws.Range("E44").Value = "  -6.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.13"
$ws.Range("E45").Value = "  -8.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.003"
$ws.Range("E46").Value = "  +0.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5973"
$ws.Range("E47").Value = "  -5.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.761"
$ws.Range("E48").Value = "  -3.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.104"
$ws.Range("E49").Value = "  -5.40%  "
$ws.Range("E50").Value = "  +3.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "124.87"
$ws.Range("E51").Value = "  -1.15%  "
